# Refresh the cryptos price/volume snapshot (Price = column D, Volume(1h) = column E).
# NumberFormat is forced to Text ("@") before writing any Price value that would
# otherwise be auto-parsed as a number (stripping trailing zeros / turning into a
# float) by Excel, so the literal display string from the source feed is preserved.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.118.71"
$ws.Range("E2").Value = "  +1.19%  "
$ws.Range("D3").Value = "3.555.00"
$ws.Range("E3").Value = "  +5.22%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.62"
$ws.Range("E5").Value = "  +1.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.27"
$ws.Range("E6").Value = "  +1.96%  "
$ws.Range("D7").Value = "3.554.21"
$ws.Range("E7").Value = "  +5.25%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.488"
$ws.Range("E9").Value = "  +4.32%  "
$ws.Range("E10").Value = "  +2.52%  "
$ws.Range("E11").Value = "  +1.03%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.411"
$ws.Range("E12").Value = "  +1.97%  "
$ws.Range("D13").Value = "4.159.85"
$ws.Range("E13").Value = "  +5.23%  "
$ws.Range("E14").Value = "  +3.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "30.00"
$ws.Range("E15").Value = "  +2.01%  "
$ws.Range("D16").Value = "3.557.76"
$ws.Range("E16").Value = "  +5.27%  "
$ws.Range("D17").Value = "66.207.94"
$ws.Range("E17").Value = "  +1.28%  "
$ws.Range("E18").Value = "  -0.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.35"
$ws.Range("E19").Value = "  +10.29%  "
$ws.Range("E20").Value = "  +1.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.85"
$ws.Range("E21").Value = "  +2.35%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "429.78"
$ws.Range("E22").Value = "  +4.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.609"
$ws.Range("E23").Value = "  +5.64%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.78"
$ws.Range("E24").Value = "  +2.56%  "
$ws.Range("D25").Value = "3.696.97"
$ws.Range("E25").Value = "  +4.97%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E27").Value = "  +7.89%  "
$ws.Range("E28").Value = "  +4.88%  "
$ws.Range("E29").Value = "  +2.56%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.09"
$ws.Range("E30").Value = "  -0.96%  "
$ws.Range("E31").Value = "  -0.09%  "
$ws.Range("E32").Value = "  +1.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "25.49"
$ws.Range("E33").Value = "  +5.03%  "
$ws.Range("D34").Value = "3.550.60"
$ws.Range("E34").Value = "  +5.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.153"
$ws.Range("E35").Value = "  -4.44%  "
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("E37").Value = "  +4.39%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.88"
$ws.Range("E38").Value = "  +5.47%  "
$ws.Range("E39").Value = "  +1.79%  "
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "170.54"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0851"
$ws.Range("E42").Value = "  +0.12%  "
$ws.Range("E43").Value = "  +3.52%  "
$ws.Range("E44").Value = "  +3.57%  "
$ws.Range("E45").Value = "  +0.44%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "46.13"
$ws.Range("E46").Value = "  +1.86%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "25.85"
$ws.Range("E47").Value = "  -1.93%  "
$ws.Range("E48").Value = "  +2.48%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.35"
$ws.Range("E49").Value = "  +4.74%  "
$ws.Range("E50").Value = "  +1.51%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.33"
$ws.Range("E51").Value = "  +15.73%  "
